$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 141.78572
$ws.Range("I4").Value = 98.833336
$ws.Range("J4").Value = 399.5
$ws.Range("K4").Value = 98.833336
$ws.Range("L4").Value = 399.5
$ws.Range("M4").Value = 15.166664
$ws.Range("N4").Value = -627.5

$ws.Range("H18").Value = 462.5
$ws.Range("I18").Value = 380
$ws.Range("K18").Value = 380
$ws.Range("M18").Value = -96

$ws.Range("H32").Value = 265.35715
$ws.Range("J32").Value = 244.85715
$ws.Range("L32").Value = 244.85715
$ws.Range("N32").Value = -896.85715

$ws.Range("H70").Value = 5545.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 5545.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 16636.5
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -17176.5

$ws.Range("H73").Value = 5545.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 5545.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 16636.5
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -18508.5

$ws.Range("H82").Value = 6329.3076
$ws.Range("I82").Value = 1840.5714
$ws.Range("J82").Value = 11566.167
$ws.Range("K82").Value = 5521.7142
$ws.Range("L82").Value = 34698.501
$ws.Range("M82").Value = -5115.7142
$ws.Range("N82").Value = -35510.501

$ws.Range("H85").Value = 6329.3076
$ws.Range("I85").Value = 1840.5714
$ws.Range("J85").Value = 11566.167
$ws.Range("K85").Value = 5521.7142
$ws.Range("L85").Value = 34698.501
$ws.Range("M85").Value = -4117.7142
$ws.Range("N85").Value = -37506.501

$ws.Range("H98").Value = 2068.4707
$ws.Range("I98").Value = 1507.4103
$ws.Range("J98").Value = 3891.9167
$ws.Range("K98").Value = 1507.4103
$ws.Range("L98").Value = 3891.9167
$ws.Range("M98").Value = -9.410300000000007
$ws.Range("N98").Value = -6887.9167

$ws.Range("H103").Value = 537.75
$ws.Range("I103").Value = 554.8182
$ws.Range("J103").Value = 350
$ws.Range("K103").Value = 1664.4546
$ws.Range("L103").Value = 1050
$ws.Range("M103").Value = -1078.4546
$ws.Range("N103").Value = -2222

$ws.Range("H122").Value = 2068.4707
$ws.Range("I122").Value = 1507.4103
$ws.Range("J122").Value = 3891.9167
$ws.Range("K122").Value = 4522.2309
$ws.Range("L122").Value = 11675.7501
$ws.Range("M122").Value = -2072.2309
$ws.Range("N122").Value = -16575.7501

$ws.Range("H132").Value = 103281.55
$ws.Range("I132").Value = 144559.53
$ws.Range("J132").Value = 6966.25
$ws.Range("K132").Value = 433678.59
$ws.Range("L132").Value = 20898.75
$ws.Range("M132").Value = -431148.59
$ws.Range("N132").Value = -25958.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4579.3
$ws.Range("I32").Value = 2881.7576
$ws.Range("J32").Value = 7874.5293
$ws.Range("K32").Value = 2881.7576
$ws.Range("L32").Value = 7874.5293
$ws.Range("M32").Value = -2594.7576
$ws.Range("N32").Value = -8448.5293

$ws.Range("H61").Value = 1985
$ws.Range("I61").Value = 1278.1177
$ws.Range("J61").Value = 3701.7144
$ws.Range("K61").Value = 1278.1177
$ws.Range("L61").Value = 3701.7144
$ws.Range("M61").Value = -1066.1177
$ws.Range("N61").Value = -4125.7144

$ws.Range("H76").Value = 19400
$ws.Range("I76").Value = 2000
$ws.Range("J76").Value = 36800
$ws.Range("K76").Value = 2000
$ws.Range("L76").Value = 36800
$ws.Range("M76").Value = -1662
$ws.Range("N76").Value = -37476

$ws.Range("H79").Value = 19400
$ws.Range("I79").Value = 2000
$ws.Range("J79").Value = 36800
$ws.Range("K79").Value = 2000
$ws.Range("L79").Value = 36800
$ws.Range("M79").Value = -830
$ws.Range("N79").Value = -39140

$ws.Range("H122").Value = 3099.2727
$ws.Range("I122").Value = 1386.5
$ws.Range("K122").Value = 4159.5
$ws.Range("M122").Value = -1709.5

$ws.Range("H132").Value = 2532.5454
$ws.Range("I132").Value = 1761.6744
$ws.Range("K132").Value = 5285.023200000001
$ws.Range("M132").Value = -2755.023200000001

$ws.Range("H136").Value = 1985
$ws.Range("I136").Value = 1278.1177
$ws.Range("J136").Value = 3701.7144
$ws.Range("K136").Value = 3834.3531
$ws.Range("L136").Value = 11105.1432
$ws.Range("M136").Value = -1284.3531
$ws.Range("N136").Value = -16205.1432

$ws.Range("H137").Value = 52780
$ws.Range("J137").Value = 52780
$ws.Range("L137").Value = 52780
$ws.Range("N137").Value = -62980

$ws.Range("H139").Value = 42869.547
$ws.Range("J139").Value = 42869.547
$ws.Range("L139").Value = 42869.547
$ws.Range("N139").Value = -53149.547

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 274.52942
$ws.Range("I22").Value = 274.52942
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 274.52942
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -101.52942
$ws.Range("N22").ClearContents()

$ws.Range("H86").Value = 2225.5715
$ws.Range("I86").Value = 2429.8333
$ws.Range("J86").Value = 1000
$ws.Range("K86").Value = 2429.8333
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = -1306.8333
$ws.Range("N86").Value = -3246

$ws.Range("H89").Value = 2225.5715
$ws.Range("I89").Value = 2429.8333
$ws.Range("J89").Value = 1000
$ws.Range("K89").Value = 12149.1665
$ws.Range("L89").Value = 5000
$ws.Range("M89").Value = -6533.166499999999
$ws.Range("N89").Value = -16232

$ws.Range("H94").Value = 1527.1305
$ws.Range("I94").Value = 1459.95
$ws.Range("J94").Value = 1975
$ws.Range("K94").Value = 1459.95
$ws.Range("L94").Value = 1975
$ws.Range("M94").Value = -1008.95
$ws.Range("N94").Value = -2877

$ws.Range("H134").Value = 2863.2354
$ws.Range("I134").Value = 1551.3
$ws.Range("J134").Value = 7633.909
$ws.Range("K134").Value = 4653.9
$ws.Range("L134").Value = 22901.727
$ws.Range("M134").Value = -2118.9
$ws.Range("N134").Value = -27971.727

$ws.Range("H138").Value = 41330.285
$ws.Range("J138").Value = 41330.285
$ws.Range("L138").Value = 41330.285
$ws.Range("N138").Value = -51610.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2505.077
$ws.Range("I31").Value = 881.15625
$ws.Range("J31").Value = 5103.35
$ws.Range("K31").Value = 881.15625
$ws.Range("L31").Value = 5103.35
$ws.Range("M31").Value = -586.15625
$ws.Range("N31").Value = -5693.35

$ws.Range("H34").Value = 2505.077
$ws.Range("I34").Value = 881.15625
$ws.Range("J34").Value = 5103.35
$ws.Range("K34").Value = 881.15625
$ws.Range("L34").Value = 5103.35
$ws.Range("M34").Value = -679.15625
$ws.Range("N34").Value = -5507.35

$ws.Range("H58").Value = 2000.6268
$ws.Range("I58").Value = 1690.1356
$ws.Range("J58").Value = 4290.5
$ws.Range("K58").Value = 1690.1356
$ws.Range("L58").Value = 4290.5
$ws.Range("M58").Value = -1487.1356
$ws.Range("N58").Value = -4696.5

$ws.Range("H111").Value = 41800
$ws.Range("J111").Value = 41800
$ws.Range("L111").Value = 41800
$ws.Range("N111").Value = -49980

$ws.Range("H122").Value = 1978.5
$ws.Range("I122").Value = 1668.9048
$ws.Range("K122").Value = 5006.7144
$ws.Range("M122").Value = -2556.7144

$ws.Range("H132").Value = 3275.2092
$ws.Range("I132").Value = 2531.84
$ws.Range("J132").Value = 4307.6665
$ws.Range("K132").Value = 7595.52
$ws.Range("L132").Value = 12922.9995
$ws.Range("M132").Value = -5065.52
$ws.Range("N132").Value = -17982.9995

$ws.Range("H134").Value = 5185.533
$ws.Range("I134").Value = 5584.4546
$ws.Range("J134").Value = 4088.5
$ws.Range("K134").Value = 16753.3638
$ws.Range("L134").Value = 12265.5
$ws.Range("M134").Value = -14218.3638
$ws.Range("N134").Value = -17335.5

$ws.Range("H136").Value = 2000.6268
$ws.Range("I136").Value = 1690.1356
$ws.Range("J136").Value = 4290.5
$ws.Range("K136").Value = 5070.406800000001
$ws.Range("L136").Value = 12871.5
$ws.Range("M136").Value = -2520.406800000001
$ws.Range("N136").Value = -17971.5

$ws.Range("H138").Value = 38163.75
$ws.Range("J138").Value = 38163.75
$ws.Range("L138").Value = 38163.75
$ws.Range("N138").Value = -48443.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 2450.9092
$ws.Range("J117").Value = 2450.9092
$ws.Range("L117").Value = 7352.7276
$ws.Range("N117").Value = -14236.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 41152
$ws.Range("J46").Value = 41152
$ws.Range("L46").Value = 41152
$ws.Range("N46").Value = -41464

$ws.Range("H122").Value = 2360.7742
$ws.Range("I122").Value = 1859.9565
$ws.Range("J122").Value = 3800.625
$ws.Range("K122").Value = 5579.8695
$ws.Range("L122").Value = 11401.875
$ws.Range("M122").Value = -3129.8695
$ws.Range("N122").Value = -16301.875

$ws.Range("H132").Value = 2691.0908
$ws.Range("I132").Value = 1594.591
$ws.Range("K132").Value = 4783.772999999999
$ws.Range("M132").Value = -2253.772999999999

$ws.Range("H137").Value = 82190
$ws.Range("J137").Value = 82190
$ws.Range("L137").Value = 82190
$ws.Range("N137").Value = -92390

$ws.Range("H140").Value = 42608.42
$ws.Range("J140").Value = 42608.42
$ws.Range("L140").Value = 42608.42
$ws.Range("N140").Value = -52968.42

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3684.25
$ws.Range("I7").Value = 1759.3334
$ws.Range("K7").Value = 1759.3334
$ws.Range("M7").Value = -1647.3334

$ws.Range("H40").Value = 3956.0588
$ws.Range("I40").Value = 3633.5122
$ws.Range("K40").Value = 3633.5122
$ws.Range("M40").Value = -3497.5122

$ws.Range("H46").Value = 1699.9667
$ws.Range("I46").Value = 1566.5555
$ws.Range("J46").Value = 1900.0834
$ws.Range("K46").Value = 1566.5555
$ws.Range("L46").Value = 1900.0834
$ws.Range("M46").Value = -1378.5555
$ws.Range("N46").Value = -2276.0834

$ws.Range("H68").Value = 796.6486
$ws.Range("I68").Value = 735.44446
$ws.Range("K68").Value = 735.44446
$ws.Range("M68").Value = 13.55553999999995

$ws.Range("H71").Value = 796.6486
$ws.Range("I71").Value = 735.44446
$ws.Range("K71").Value = 3677.2223
$ws.Range("M71").Value = 66.77769999999964

$ws.Range("H93").Value = 2805.2144
$ws.Range("I93").Value = 1974.7778
$ws.Range("J93").Value = 4300
$ws.Range("K93").Value = 1974.7778
$ws.Range("L93").Value = 4300
$ws.Range("M93").Value = -726.7778000000001
$ws.Range("N93").Value = -6796

$ws.Range("H126").Value = 3684.25
$ws.Range("I126").Value = 1759.3334
$ws.Range("K126").Value = 5278.0002
$ws.Range("M126").Value = -2808.0002

$ws.Range("H132").Value = 5799.7744
$ws.Range("I132").Value = 2106.0527
$ws.Range("J132").Value = 11648.167
$ws.Range("K132").Value = 6318.158100000001
$ws.Range("L132").Value = 34944.501
$ws.Range("M132").Value = -3788.158100000001
$ws.Range("N132").Value = -40004.501

$ws.Range("H136").Value = 3740.9119
$ws.Range("I136").Value = 1389.55
$ws.Range("J136").Value = 7100
$ws.Range("K136").Value = 4168.65
$ws.Range("L136").Value = 21300
$ws.Range("M136").Value = -1618.65
$ws.Range("N136").Value = -26400

$ws.Range("H140").Value = 62889.75
$ws.Range("J140").Value = 62889.75
$ws.Range("L140").Value = 62889.75
$ws.Range("N140").Value = -73249.75

$ws.Range("H141").Value = 41323.055
$ws.Range("J141").Value = 41323.055
$ws.Range("L141").Value = 41323.055
$ws.Range("N141").Value = -51683.055

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 45466800
$ws.Range("I62").Value = 125003000
$ws.Range("J62").Value = 17542.857
$ws.Range("K62").Value = 125003000
$ws.Range("L62").Value = 17542.857
$ws.Range("M62").Value = -125002376
$ws.Range("N62").Value = -18790.857

$ws.Range("H65").Value = 45466800
$ws.Range("I65").Value = 125003000
$ws.Range("J65").Value = 17542.857
$ws.Range("K65").Value = 625015000
$ws.Range("L65").Value = 87714.285
$ws.Range("M65").Value = -625011880
$ws.Range("N65").Value = -93954.285

$ws.Range("H81").Value = 1583.3334
$ws.Range("I81").Value = 1875
$ws.Range("K81").Value = 3750
$ws.Range("M81").Value = -2689

$ws.Range("H84").Value = 1583.3334
$ws.Range("I84").Value = 1875
$ws.Range("K84").Value = 18750
$ws.Range("M84").Value = -13446

$ws.Range("H132").Value = 7753596
$ws.Range("I132").Value = 745.4400000000001
$ws.Range("K132").Value = 2236.32
$ws.Range("M132").Value = 293.6799999999998

$ws.Range("H136").Value = 2999.2646
$ws.Range("I136").Value = 734.65216
$ws.Range("J136").Value = 7734.364
$ws.Range("K136").Value = 2203.95648
$ws.Range("L136").Value = 23203.092
$ws.Range("M136").Value = 346.0435200000002
$ws.Range("N136").Value = -28303.092

$ws.Range("H138").Value = 52207.145
$ws.Range("J138").Value = 52207.145
$ws.Range("L138").Value = 52207.145
$ws.Range("N138").Value = -62487.145

$ws.Range("H139").Value = 39358.438
$ws.Range("J139").Value = 39627.332
$ws.Range("L139").Value = 39627.332
$ws.Range("N139").Value = -49907.332

$ws.Range("H141").Value = 43261.2
$ws.Range("J141").Value = 43261.2
$ws.Range("L141").Value = 43261.2
$ws.Range("N141").Value = -53621.2

Write-Host "Updated Chocobo Profits figures across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets"
